$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.6183950451608057
$ws.Range("E2").Value = 0.0007819828803196954
$ws.Range("F2").Value = 24.96723176342669
$ws.Range("I2").Value = 1.611381718659552
$ws.Range("J2").Value = "[ 0.31778082 -0.06268567]"
$ws.Range("D3").Value = 0.6325629441036188
$ws.Range("E3").Value = 0.01836892643032831
$ws.Range("F3").Value = 16.73684223409452
$ws.Range("I3").Value = 1.080193505617841
$ws.Range("J3").Value = "[0.176319   0.04110187]"
$ws.Range("D4").Value = 0.7024686972505991
$ws.Range("E4").Value = 0.3814594999279021
$ws.Range("F4").Value = 9.977897272211113
$ws.Range("I4").Value = 0.6439721234396593
$ws.Range("J4").Value = "[0.02421953 0.12041552]"
$ws.Range("D5").Value = 0.6035278586660175
$ws.Range("E5").Value = 0.0006166715590086611
$ws.Range("F5").Value = 20.22047146842954
$ws.Range("I5").Value = 1.305026459306296
$ws.Range("J5").Value = "[ 0.19538445 -0.25052374]"
$ws.Range("D6").Value = 0.675418437646664
$ws.Range("E6").Value = 0.002798440864392866
$ws.Range("F6").Value = 18.89071872877038
$ws.Range("I6").Value = 1.219204399603104
$ws.Range("J6").Value = "[ 0.13661115 -0.3502805 ]"
$ws.Range("D7").Value = 0.8712663961374646
$ws.Range("E7").Value = 0.001267772930683297
$ws.Range("F7").Value = 18.6059829129727
$ws.Range("I7").Value = 1.200827589047114
$ws.Range("J7").Value = "[ 0.43588112 -0.54375427]"
$ws.Range("D8").Value = 0.8480506102993624
$ws.Range("E8").Value = 0.01197733081141318
$ws.Range("F8").Value = 16.0055076266447
$ws.Range("I8").Value = 1.03299327021193
$ws.Range("J8").Value = "[ 0.37952504 -0.46381184]"
$ws.Range("D9").Value = 0.7302709731347655
$ws.Range("E9").Value = 0.00954509154990923
$ws.Range("F9").Value = 15.67207963345827
$ws.Range("I9").Value = 1.011473873195838
$ws.Range("J9").Value = "[ 0.37999036 -0.04464736]"
$ws.Range("D10").Value = 0.7964450256248707
$ws.Range("E10").Value = 0.02341053620396535
$ws.Range("F10").Value = 10.45988880418429
$ws.Range("I10").Value = 0.6750797909027378
$ws.Range("J10").Value = "[ 0.40515241 -0.3874311 ]"
$ws.Range("D11").Value = 0.7386033464125396
$ws.Range("E11").Value = 0.01826050842569454
$ws.Range("F11").Value = 12.30909841830705
$ws.Range("I11").Value = 0.7944275261423266
$ws.Range("J11").Value = "[ 0.41547965 -0.03906522]"
$ws.Range("D12").Value = 0.7782430363992037
$ws.Range("E12").Value = 0.0295349030511869
$ws.Range("F12").Value = 10.78444352841938
$ws.Range("I12").Value = 0.6960265083559363
$ws.Range("J12").Value = "[ 0.33661073 -0.03340498]"
$ws.Range("D13").Value = 0.7515369188484933
$ws.Range("E13").Value = 0.01912841741346826
$ws.Range("F13").Value = 18.37780103456878
$ws.Range("I13").Value = 1.186100761865252
$ws.Range("J13").Value = "[ 0.21263165 -0.31504723]"
$ws.Range("D14").Value = 0.772185228876031
$ws.Range("E14").Value = 0.03016498529027952
$ws.Range("F14").Value = 9.293349210705664
$ws.Range("I14").Value = 0.5997914853014117
$ws.Range("J14").Value = "[ 0.26785419 -0.09700592]"
$ws.Range("D15").Value = 0.5505465446083191
$ws.Range("E15").Value = 0.06868595729533095
$ws.Range("F15").Value = 7.381518027880717
$ws.Range("I15").Value = 0.4764021625940325
$ws.Range("J15").Value = "[ 0.21649303 -0.02877684]"
$ws.Range("D16").Value = 0.6020397476619385
$ws.Range("E16").Value = 0.06314370940112339
$ws.Range("F16").Value = 8.524170068501276
$ws.Range("I16").Value = 0.5501487688053849
$ws.Range("D17").Value = 0.3725246591215406
$ws.Range("E17").Value = 0.08313840402119654
$ws.Range("F17").Value = 7.319725145562557
$ws.Range("I17").Value = 0.4724140584319754
$ws.Range("J17").Value = "[ 0.19751082 -0.01785914]"
$ws.Range("D18").Value = 0.7658831334318821
$ws.Range("E18").Value = 0.01911849081782102
$ws.Range("F18").Value = 15.74098723456522
$ws.Range("I18").Value = 1.015921160334137
$ws.Range("J18").Value = "[ 0.18011774 -0.06764064]"
$ws.Range("D19").Value = 0.6021386100779408
$ws.Range("E19").Value = 0.0008485706930041688
$ws.Range("F19").Value = 18.15952315501167
$ws.Range("I19").Value = 1.172013137412576
$ws.Range("J19").Value = "[ 0.55961215 -0.06370079]"
$ws.Range("D20").Value = 0.6099209620270115
$ws.Range("E20").Value = 0.00202742775725777
$ws.Range("F20").Value = 17.59093869546037
$ws.Range("I20").Value = 1.135316774262817
$ws.Range("J20").Value = "[ 0.47875293 -0.20064185]"
$ws.Range("D21").Value = 0.4275101007180105
$ws.Range("E21").Value = 0.03383070218535345
$ws.Range("F21").Value = 10.48518143524562
$ws.Range("I21").Value = 0.6767121738475191
$ws.Range("J21").Value = "[ 0.36061115 -0.50433351]"
$ws.Range("D22").Value = 0.5978223617730578
$ws.Range("E22").Value = 0.01126703751299983
$ws.Range("F22").Value = 12.38659587350982
$ws.Range("I22").Value = 0.7994292012875585
$ws.Range("J22").Value = "[ 0.4710036  -0.18720615]"
$ws.Range("D23").Value = 0.6269290030932835
$ws.Range("E23").Value = 0.2750046805800044
$ws.Range("F23").Value = 7.097662131966908
$ws.Range("I23").Value = 0.458082141946839
$ws.Range("J23").Value = "[ 0.12747955 -0.40589612]"
$ws.Range("D24").Value = 0.6420917296715405
$ws.Range("E24").Value = 0.1013514275555494
$ws.Range("F24").Value = 7.818924508511154
$ws.Range("I24").Value = 0.5046323169495333
$ws.Range("J24").Value = "[ 0.1408866  -0.12248282]"
$ws.Range("D25").Value = 0.6850899204272667
$ws.Range("E25").Value = 0.001150516603473799
$ws.Range("F25").Value = 16.78219035559083
$ws.Range("I25").Value = 1.083120267168626
$ws.Range("J25").Value = "[ 0.76322768 -0.22079373]"
$ws.Range("D26").Value = 0.7037610421524042
$ws.Range("E26").Value = 0.00248283298082216
$ws.Range("F26").Value = 17.44924399721732
$ws.Range("I26").Value = 1.126171817843808
$ws.Range("J26").Value = "[ 0.60993771 -0.71030696]"
$ws.Range("D27").Value = 0.6463059737201324
$ws.Range("E27").Value = 0.003636994000681314
$ws.Range("F27").Value = 15.61975608935025
$ws.Range("I27").Value = 1.008096918825003
